{"js": "// Replace each two-digit multiplication equation in the document's table\n// with its updated equation/result, per the commit diff.\nconst replacements = [\n  [\"40\u00d737=1480\", \"46\u00d758=2668\"],\n  [\"37\u00d724=888\", \"74\u00d734=2516\"],\n  [\"61\u00d727=1647\", \"29\u00d775=2175\"],\n  [\"92\u00d731=2852\", \"77\u00d737=2849\"],\n  [\"56\u00d784=4704\", \"11\u00d778=858\"],\n  [\"35\u00d728=980\", \"76\u00d741=3116\"],\n  [\"49\u00d782=4018\", \"66\u00d754=3564\"],\n  [\"24\u00d713=312\", \"48\u00d737=1776\"],\n  [\"79\u00d714=1106\", \"92\u00d797=8924\"],\n  [\"66\u00d758=3828\", \"25\u00d777=1925\"],\n  [\"24\u00d789=2136\", \"65\u00d729=1885\"],\n  [\"17\u00d751=867\", \"82\u00d716=1312\"],\n  [\"67\u00d759=3953\", \"64\u00d749=3136\"],\n  [\"90\u00d723=2070\", \"96\u00d735=3360\"],\n  [\"12\u00d715=180\", \"11\u00d746=506\"],\n  [\"16\u00d791=1456\", \"88\u00d797=8536\"],\n  [\"54\u00d796=5184\", \"96\u00d718=1728\"],\n  [\"43\u00d785=3655\", \"59\u00d782=4838\"],\n  [\"43\u00d754=2322\", \"97\u00d742=4074\"],\n  [\"30\u00d744=1320\", \"53\u00d732=1696\"],\n  [\"88\u00d775=6600\", \"69\u00d770=4830\"],\n  [\"95\u00d744=4180\", \"26\u00d782=2132\"],\n  [\"77\u00d794=7238\", \"94\u00d746=4324\"],\n  [\"18\u00d711=198\", \"50\u00d791=4550\"],\n  [\"76\u00d737=2812\", \"95\u00d730=2850\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document's table\n# with its updated equation/result, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"40\u00d737=1480\"; New = \"46\u00d758=2668\" },\n    @{ Old = \"37\u00d724=888\";  New = \"74\u00d734=2516\" },\n    @{ Old = \"61\u00d727=1647\"; New = \"29\u00d775=2175\" },\n    @{ Old = \"92\u00d731=2852\"; New = \"77\u00d737=2849\" },\n    @{ Old = \"56\u00d784=4704\"; New = \"11\u00d778=858\" },\n    @{ Old = \"35\u00d728=980\";  New = \"76\u00d741=3116\" },\n    @{ Old = \"49\u00d782=4018\"; New = \"66\u00d754=3564\" },\n    @{ Old = \"24\u00d713=312\";  New = \"48\u00d737=1776\" },\n    @{ Old = \"79\u00d714=1106\"; New = \"92\u00d797=8924\" },\n    @{ Old = \"66\u00d758=3828\"; New = \"25\u00d777=1925\" },\n    @{ Old = \"24\u00d789=2136\"; New = \"65\u00d729=1885\" },\n    @{ Old = \"17\u00d751=867\";  New = \"82\u00d716=1312\" },\n    @{ Old = \"67\u00d759=3953\"; New = \"64\u00d749=3136\" },\n    @{ Old = \"90\u00d723=2070\"; New = \"96\u00d735=3360\" },\n    @{ Old = \"12\u00d715=180\";  New = \"11\u00d746=506\" },\n    @{ Old = \"16\u00d791=1456\"; New = \"88\u00d797=8536\" },\n    @{ Old = \"54\u00d796=5184\"; New = \"96\u00d718=1728\" },\n    @{ Old = \"43\u00d785=3655\"; New = \"59\u00d782=4838\" },\n    @{ Old = \"43\u00d754=2322\"; New = \"97\u00d742=4074\" },\n    @{ Old = \"30\u00d744=1320\"; New = \"53\u00d732=1696\" },\n    @{ Old = \"88\u00d775=6600\"; New = \"69\u00d770=4830\" },\n    @{ Old = \"95\u00d744=4180\"; New = \"26\u00d782=2132\" },\n    @{ Old = \"77\u00d794=7238\"; New = \"94\u00d746=4324\" },\n    @{ Old = \"18\u00d711=198\";  New = \"50\u00d791=4550\" },\n    @{ Old = \"76\u00d737=2812\"; New = \"95\u00d730=2850\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
